# "Added more wind data"
# Adds 7 new weather stations (Leeming, Linton, Waddington, Wainfleet,
# Barkston, Cranwell, Wittering) to the bottom of the wind-speed table,
# fills in the previously-missing P64 average, extends the P-column
# AVERAGE() formulas down through the new rows, and leaves the selection
# on Q66 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64 already had its 12 monthly values but was missing its Average
# (column P) formula - the diff backfills it.
$ws.Range("P64").Formula = "=AVERAGE(D64:O64)"

# --- Leeming ---------------------------------------------------------
$ws.Range("A65").Value = "Leeming"
$ws.Range("B65").Value = 54.3
$ws.Range("C65").Value = -1.5333000000000001
$ws.Range("D65").Value = 5.0999999999999996
$ws.Range("E65").Value = 3.94
$ws.Range("F65").Value = 4.97
$ws.Range("G65").Value = 3.88
$ws.Range("H65").Value = 4.07
$ws.Range("I65").Value = 3.34
$ws.Range("J65").Value = 3.68
$ws.Range("K65").Value = 3.86
$ws.Range("L65").Value = 3.86
$ws.Range("M65").Value = 3.73
$ws.Range("N65").Value = 4.78
$ws.Range("O65").Value = 3.74
$ws.Range("P65").Formula = "=AVERAGE(D65:O65)"

# --- Linton ------------------------------------------------------------
$ws.Range("A66").Value = "Linton"
$ws.Range("B66").Value = 54.05
$ws.Range("C66").Value = -1.25
$ws.Range("D66").Value = 4.3600000000000003
$ws.Range("E66").Value = 3.52
$ws.Range("F66").Value = 4.5199999999999996
$ws.Range("G66").Value = 3.6
$ws.Range("H66").Value = 3.79
$ws.Range("I66").Value = 2.97
$ws.Range("J66").Value = 3.42
$ws.Range("K66").Value = 3.52
$ws.Range("L66").Value = 3.44
$ws.Range("M66").Value = 3.46
$ws.Range("N66").Value = 4.3899999999999997
$ws.Range("O66").Value = 3.38
$ws.Range("P66").Formula = "=AVERAGE(D66:O66)"

# --- Waddington ----------------------------------------------------------
$ws.Range("A67").Value = "Waddington"
$ws.Range("B67").Value = 53.166699999999999
$ws.Range("C67").Value = -0.51670000000000005
$ws.Range("D67").Value = 5.3
$ws.Range("E67").Value = 4.2699999999999996
$ws.Range("F67").Value = 5.0199999999999996
$ws.Range("G67").Value = 4.42
$ws.Range("H67").Value = 4.55
$ws.Range("I67").Value = 3.67
$ws.Range("J67").Value = 4.22
$ws.Range("K67").Value = 4.18
$ws.Range("L67").Value = 4.28
$ws.Range("M67").Value = 4.41
$ws.Range("N67").Value = 5.19
$ws.Range("O67").Value = 4.45
$ws.Range("P67").Formula = "=AVERAGE(D67:O67)"

# --- Wainfleet -----------------------------------------------------------
$ws.Range("A68").Value = "Wainfleet"
$ws.Range("B68").Value = 53.083300000000001
$ws.Range("C68").Value = 0.26669999999999999
$ws.Range("D68").Value = 6.51
$ws.Range("E68").Value = 5.34
$ws.Range("F68").Value = 6.25
$ws.Range("G68").Value = 5.54
$ws.Range("H68").Value = 5.85
$ws.Range("I68").Value = 4.54
$ws.Range("J68").Value = 4.87
$ws.Range("K68").Value = 4.95
$ws.Range("L68").Value = 4.63
$ws.Range("M68").Value = 5.48
$ws.Range("N68").Value = 6.45
$ws.Range("O68").Value = 5.86
$ws.Range("P68").Formula = "=AVERAGE(D68:O68)"

# --- Barkston ------------------------------------------------------------
$ws.Range("A69").Value = "Barkston"
$ws.Range("B69").Value = 52.966700000000003
$ws.Range("C69").Value = -0.56669000000000003
$ws.Range("D69").Value = 7.08
$ws.Range("E69").Value = 5.83
$ws.Range("F69").Value = 6.76
$ws.Range("G69").Value = 5.97
$ws.Range("H69").Value = 6.06
$ws.Range("I69").Value = 4.63
$ws.Range("J69").Value = 5.27
$ws.Range("K69").Value = 5.61
$ws.Range("L69").Value = 5.81
$ws.Range("M69").Value = 5.48
$ws.Range("N69").Value = 6.67
$ws.Range("O69").Value = 6.19
$ws.Range("P69").Formula = "=AVERAGE(D69:O69)"

# --- Cranwell ------------------------------------------------------------
$ws.Range("A70").Value = "Cranwell"
$ws.Range("B70").Value = 53.033000000000001
$ws.Range("C70").Value = -0.5
$ws.Range("D70").Value = 6.23
$ws.Range("E70").Value = 4.96
$ws.Range("F70").Value = 6.02
$ws.Range("G70").Value = 4.8600000000000003
$ws.Range("H70").Value = 4.92
$ws.Range("I70").Value = 3.97
$ws.Range("J70").Value = 4.79
$ws.Range("K70").Value = 4.93
$ws.Range("L70").Value = 4.82
$ws.Range("M70").Value = 4.87
$ws.Range("N70").Value = 6.1
$ws.Range("O70").Value = 5.35
$ws.Range("P70").Formula = "=AVERAGE(D70:O70)"

# --- Wittering -----------------------------------------------------------
$ws.Range("A71").Value = "Wittering"
$ws.Range("B71").Value = 52.616700000000002
$ws.Range("C71").Value = -0.46668999999999999
$ws.Range("D71").Value = 5.8
$ws.Range("E71").Value = 4.74
$ws.Range("F71").Value = 5.78
$ws.Range("G71").Value = 4.84
$ws.Range("H71").Value = 5.0199999999999996
$ws.Range("I71").Value = 4.09
$ws.Range("J71").Value = 4.71
$ws.Range("K71").Value = 4.7699999999999996
$ws.Range("L71").Value = 4.62
$ws.Range("M71").Value = 4.57
$ws.Range("N71").Value = 5.58
$ws.Range("O71").Value = 4.97
$ws.Range("P71").Formula = "=AVERAGE(D71:O71)"

# Conditional formatting ranges that previously stopped at row 63 grow to
# cover the new rows through 71 (same rules, just a bigger sqref).
$fc = $ws.Range("P3:P63").FormatConditions
$fc.Item(1).ModifyAppliesToRange($ws.Range("P3:P71"))
$fc.Item(2).ModifyAppliesToRange($ws.Range("P17:P71"))
$fc.Item(3).ModifyAppliesToRange($ws.Range("P17:P71"))

# Leave the cursor where the author left it.
$ws.Range("Q66").Select()
